{"js": "// Replace each old text value with its corresponding new text value.\n// Each pair is a unique, literal (non-wildcard) string in the document, so a\n// simple body.search + Replace works without touching unrelated runs.\nconst replacements = [\n  [\"2025-05-06 Tuesday\", \"2025-05-07 Wednesday\"],\n  [\"215\u00f75=43, 0\", \"638\u00f79=70, 8\"],\n  [\"632\u00f75=126, 2\", \"282\u00f76=47, 0\"],\n  [\"889\u00f77=127, 0\", \"615\u00f75=123, 0\"],\n  [\"376\u00f78=47, 0\", \"123\u00f74=30, 3\"],\n  [\"120\u00f75=24, 0\", \"537\u00f78=67, 1\"],\n  [\"899\u00f73=299, 2\", \"765\u00f76=127, 3\"],\n  [\"440\u00f73=146, 2\", \"723\u00f74=180, 3\"],\n  [\"844\u00f78=105, 4\", \"821\u00f75=164, 1\"],\n  [\"389\u00f74=97, 1\", \"150\u00f75=30, 0\"],\n  [\"361\u00f76=60, 1\", \"713\u00f73=237, 2\"],\n  [\"255\u00f76=42, 3\", \"564\u00f72=282, 0\"],\n  [\"620\u00f73=206, 2\", \"943\u00f79=104, 7\"],\n  [\"497\u00f79=55, 2\", \"486\u00f74=121, 2\"],\n  [\"628\u00f77=89, 5\", \"710\u00f79=78, 8\"],\n  [\"236\u00f72=118, 0\", \"795\u00f73=265, 0\"],\n  [\"865\u00f74=216, 1\", \"689\u00f72=344, 1\"],\n  [\"823\u00f73=274, 1\", \"791\u00f72=395, 1\"],\n  [\"166\u00f76=27, 4\", \"632\u00f77=90, 2\"],\n  [\"525\u00f76=87, 3\", \"897\u00f74=224, 1\"],\n  [\"676\u00f73=225, 1\", \"839\u00f73=279, 2\"],\n  [\"807\u00f76=134, 3\", \"922\u00f77=131, 5\"],\n  [\"997\u00f76=166, 1\", \"791\u00f75=158, 1\"],\n  [\"804\u00f74=201, 0\", \"143\u00f75=28, 3\"],\n  [\"262\u00f76=43, 4\", \"335\u00f78=41, 7\"],\n  [\"211\u00f77=30, 1\", \"726\u00f72=363, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n", "ps1": "# Replace each old three-digit-division answer (and the date line) with\n# its corresponding new value. Every \"old\" string is a unique, literal\n# run of text in the document, so Find/Replace (wdReplaceAll, no wildcards)\n# swaps exactly one run each and leaves all formatting untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-05-06 Tuesday\", \"2025-05-07 Wednesday\"),\n    @(\"215\u00f75=43, 0\", \"638\u00f79=70, 8\"),\n    @(\"632\u00f75=126, 2\", \"282\u00f76=47, 0\"),\n    @(\"889\u00f77=127, 0\", \"615\u00f75=123, 0\"),\n    @(\"376\u00f78=47, 0\", \"123\u00f74=30, 3\"),\n    @(\"120\u00f75=24, 0\", \"537\u00f78=67, 1\"),\n    @(\"899\u00f73=299, 2\", \"765\u00f76=127, 3\"),\n    @(\"440\u00f73=146, 2\", \"723\u00f74=180, 3\"),\n    @(\"844\u00f78=105, 4\", \"821\u00f75=164, 1\"),\n    @(\"389\u00f74=97, 1\", \"150\u00f75=30, 0\"),\n    @(\"361\u00f76=60, 1\", \"713\u00f73=237, 2\"),\n    @(\"255\u00f76=42, 3\", \"564\u00f72=282, 0\"),\n    @(\"620\u00f73=206, 2\", \"943\u00f79=104, 7\"),\n    @(\"497\u00f79=55, 2\", \"486\u00f74=121, 2\"),\n    @(\"628\u00f77=89, 5\", \"710\u00f79=78, 8\"),\n    @(\"236\u00f72=118, 0\", \"795\u00f73=265, 0\"),\n    @(\"865\u00f74=216, 1\", \"689\u00f72=344, 1\"),\n    @(\"823\u00f73=274, 1\", \"791\u00f72=395, 1\"),\n    @(\"166\u00f76=27, 4\", \"632\u00f77=90, 2\"),\n    @(\"525\u00f76=87, 3\", \"897\u00f74=224, 1\"),\n    @(\"676\u00f73=225, 1\", \"839\u00f73=279, 2\"),\n    @(\"807\u00f76=134, 3\", \"922\u00f77=131, 5\"),\n    @(\"997\u00f76=166, 1\", \"791\u00f75=158, 1\"),\n    @(\"804\u00f74=201, 0\", \"143\u00f75=28, 3\"),\n    @(\"262\u00f76=43, 4\", \"335\u00f78=41, 7\"),\n    @(\"211\u00f77=30, 1\", \"726\u00f72=363, 0\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $oldText,   # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n}\n"}
